# Scheduled Sophia-Profits refresh: update per-leve currentAveragePrice /
# LevePrice / LeveProfit figures (cols H:N) across all job sheets to match
# the latest market-board pull. Identity columns A:G are untouched.
$wb = $excel.ActiveWorkbook

# ALC!42 - "Eye of the Beholder" (Hi-Potion of Dexterity)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value = 331.83334
$ws.Range("I42").Value = 205.33333
$ws.Range("J42").Value = 458.33334
$ws.Range("K42").Value = 615.99999
$ws.Range("L42").Value = 1375.00002
$ws.Range("M42").Value = -385.99999
$ws.Range("N42").Value = -1835.00002

# ALC!51 - "A Bile Business" (Shark Oil)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 1500
$ws.Range("J51").Value = 1500
$ws.Range("L51").Value = 1500
$ws.Range("N51").Value = -2468

# ALC!92 - "Whinier than the Sword" (Enchanted Koppranickel Ink)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 786.7273
$ws.Range("I92").Value = 890.125
$ws.Range("J92").Value = 511
$ws.Range("K92").Value = 890.125
$ws.Range("L92").Value = 511
$ws.Range("M92").Value = 357.875
$ws.Range("N92").Value = -3007

# ALC!96 - "Scroll Down" (Grade 1 Reisui of Intelligence)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 421.92307
$ws.Range("I96").Value = 296.81818
$ws.Range("J96").Value = 1110
$ws.Range("K96").Value = 890.45454
$ws.Range("L96").Value = 3330
$ws.Range("M96").Value = 482.54546
$ws.Range("N96").Value = -6076

# ALC!132 - "Fast-forwarding Flora" (Growth Formula Lambda)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1289.4762
$ws.Range("I132").Value = 1289.4762
$ws.Range("K132").Value = 3868.4286
$ws.Range("M132").Value = -1338.4286

# ARM!32 - "Ingot We Trust" (Steel Ingot)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7744.033
$ws.Range("I32").Value = 7238.25
$ws.Range("K32").Value = 7238.25
$ws.Range("M32").Value = -6951.25

# ARM!50 - "Liquid Persuasion" (Mythril Chain Coif)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H50").Value = 6750
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 6750
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 6750
$ws.Range("M50").ClearContents()
$ws.Range("N50").Value = -8178

# ARM!61 - "Dealing with the Tough Stuff" (Cobalt Ingot)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2284.5715
$ws.Range("I61").Value = 2165.3333
$ws.Range("K61").Value = 2165.3333
$ws.Range("M61").Value = -1953.3333

# ARM!110 - "Scheduled Maintenance" (Deepgold Ingot)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 3174.6667
$ws.Range("I110").Value = 1179.6
$ws.Range("K110").Value = 1179.6
$ws.Range("M110").Value = 865.4000000000001

# ARM!136 - "Metal with Mettle" (Cobalt Tungsten Ingot)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2284.5715
$ws.Range("I136").Value = 2165.3333
$ws.Range("K136").Value = 6495.999899999999
$ws.Range("M136").Value = -3945.999899999999

# BSM!134 - "Ruthenium Supremium" (Ruthenium Ingot)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2487.9473
$ws.Range("I134").Value = 2348.389
$ws.Range("K134").Value = 7045.167
$ws.Range("M134").Value = -4510.167

# CRP!12 - "A Sword in Hand" (Ash Macuahuitl)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 8002.5
$ws.Range("I12").Value = 8002.5
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 8002.5
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -7832.5
$ws.Range("N12").ClearContents()

# CRP!22 - "Driving Up the Wall" (Elm Lumber)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()

# CRP!122 - "Timber of Tenkonto" (Horse Chestnut Lumber)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 4362.8184
$ws.Range("I122").Value = 3623.875
$ws.Range("J122").Value = 6333.3335
$ws.Range("K122").Value = 10871.625
$ws.Range("L122").Value = 19000.0005
$ws.Range("M122").Value = -8421.625
$ws.Range("N122").Value = -23900.0005

# CRP!134 - "Wood You Be Quiet" (Ceiba Lumber)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2940.8
$ws.Range("I134").Value = 2903.6667
$ws.Range("J134").Value = 2996.5
$ws.Range("K134").Value = 8711.000100000001
$ws.Range("L134").Value = 8989.5
$ws.Range("M134").Value = -6176.000100000001
$ws.Range("N134").Value = -14059.5

# CUL!107 - "Slippery Service" (Frantoio Oil)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 649.5
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("N107").ClearContents()

# CUL!121 - "A Cookie for Your Troubles" (Coffee Biscuit)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 1852.25
$ws.Range("J121").Value = 2106.4443
$ws.Range("L121").Value = 6319.3329
$ws.Range("N121").Value = -8939.332900000001

# CUL!132 - "More Mezcal" (Cooking Mezcal)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1527.6
$ws.Range("J132").Value = 1999.5
$ws.Range("L132").Value = 17995.5
$ws.Range("N132").Value = -23055.5

# CUL!137 - "Creative Chocolate" (Gateau au Chocolat)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 1880.625
$ws.Range("I137").Value = 1128.75
$ws.Range("J137").Value = 2632.5
$ws.Range("K137").Value = 3386.25
$ws.Range("L137").Value = 7897.5
$ws.Range("M137").Value = 1713.75
$ws.Range("N137").Value = -18097.5

# CUL!140 - "Sweet, Sweet Bean Juice" (Mesquite Juice)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 750.5
$ws.Range("I140").Value = 750.5
$ws.Range("K140").Value = 2251.5
$ws.Range("M140").Value = 2928.5

# CUL!141 - "Ocean Explosion" (Acqua Pazza)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H141").Value = 6749.25
$ws.Range("I141").Value = 6749.25
$ws.Range("K141").Value = 20247.75
$ws.Range("M141").Value = -15067.75

# GSM!70 - "Sky Is the Limit" (Mythrite Ingot)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 10000
$ws.Range("I70").Value = 10000
$ws.Range("K70").Value = 10000
$ws.Range("M70").Value = -9730

# GSM!73 - "Hulls of Broken Dreams (L)" (Mythrite Ingot)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 10000
$ws.Range("I73").Value = 10000
$ws.Range("K73").Value = 10000
$ws.Range("M73").Value = -9064

# GSM!80 - "Needs More Prayerbell" (Hardsilver Ingot)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3591.4
$ws.Range("J80").Value = 4231.5
$ws.Range("L80").Value = 4231.5
$ws.Range("N80").Value = -6227.5

# GSM!83 - "With a Noise That Reaches Heaven (L)" (Hardsilver Ingot)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 3591.4
$ws.Range("J83").Value = 4231.5
$ws.Range("L83").Value = 21157.5
$ws.Range("N83").Value = -31141.5

# GSM!113 - "Copious Crystal Cannons" (Manasilver Nugget)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1680.5714
$ws.Range("I113").Value = 1627.3334
$ws.Range("K113").Value = 1627.3334
$ws.Range("M113").Value = 542.6666

# GSM!141 - "Mask Maker" (Black Star Mask of Casting)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H141").Value = 69000
$ws.Range("J141").Value = 69000
$ws.Range("L141").Value = 69000
$ws.Range("N141").Value = -79360

# LTW!7 - "Tan Before the Ban" (Leather)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3195
$ws.Range("I7").Value = 3195
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 3195
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -3083
$ws.Range("N7").ClearContents()

# LTW!122 - "Hell on Leather" (Gaja Leather)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2999
$ws.Range("I122").Value = 2999
$ws.Range("K122").Value = 8997
$ws.Range("M122").Value = -6547

# LTW!126 - "Battered Books" (Saiga Leather)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 3195
$ws.Range("I126").Value = 3195
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 9585
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -7115
$ws.Range("N126").ClearContents()

# WVR!122 - "Heavy Armoire" (Dark Hempen Cloth)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1699.8
$ws.Range("I122").Value = 1699.8
$ws.Range("K122").Value = 5099.4
$ws.Range("M122").Value = -2649.4

# WVR!126 - "A Polished Purchase" (Snow Linen)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1630
$ws.Range("I126").Value = 1000
$ws.Range("J126").Value = 1700
$ws.Range("K126").Value = 3000
$ws.Range("L126").Value = 5100
$ws.Range("M126").Value = -530
$ws.Range("N126").Value = -10040

# WVR!132 - "Comfy Cabins" (Snow Cotton Cloth)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3325.3635
$ws.Range("I132").Value = 3271.75
$ws.Range("K132").Value = 9815.25
$ws.Range("M132").Value = -7285.25

# WVR!136 - "Weaving the Envelope" (Sarcenet Cloth)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 8998.833000000001
$ws.Range("I136").Value = 8798.6
$ws.Range("J136").Value = 10000
$ws.Range("K136").Value = 26395.8
$ws.Range("L136").Value = 30000
$ws.Range("M136").Value = -23845.8
$ws.Range("N136").Value = -35100
